$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (ErroresCarga -> Formato)
$ws.Name = "Formato"

# Force text format on numeric-looking values so leading zeros/text type are preserved
$ws.Range("A2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

# Update the data row (row 2) values
$ws.Range("A2").Value = "60000050"
$ws.Range("B2").Value = "Leche Semidescremada Tetrafino 900ml"
$ws.Range("C2").Value = "Christopher Hablich"
$ws.Range("D2").Value = "0921821419"
$ws.Range("E2").Value = "120"
$ws.Range("F2").Value = "El distribuidor no se encuentra registrado-No tiene suficiente STOCK. Cantidad actual: 42,0000"
